# Update attendance/view-count figures (column F) on the "展览" and
# "全部类型" worksheets to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 6874
    $ws.Range("F4").Value = 200
    $ws.Range("F5").Value = 35
    $ws.Range("F6").Value = 1069
    $ws.Range("F7").Value = 156
}
